$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "293.79"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.18%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.15%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07325"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-7.26%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.820"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-14.00%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.680"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.21%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.763"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.84%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9051"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1657"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-5.46%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07491"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-6.03%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08111"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-7.71%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02983"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.10%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09998"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.36%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001499"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.13%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005728"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.32%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.461"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.107"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-7.63%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.38%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1306"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.27%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.332"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.34%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "11.88%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04477"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.92%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004047"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-10.32%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001250"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01651"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.94%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04399"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-7.74%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007394"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.11%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.45%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002038"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-13.19%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01124"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.30%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006006"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.81%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.12%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "161.77%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002401"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-29.28%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.12%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"
